$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.569.93"
$ws.Range("E2").Value = "  +1.56%  "

$ws.Range("D3").Value = "2.251.45"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'246.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").Value = "'76.71"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.63%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.629"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.32%  "

$ws.Range("D10").Value = "'45.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.92%  "

$ws.Range("D11").Value = "'0.0952"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").Value = "'7.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.76%  "

$ws.Range("E13").Value = "  -0.75%  "

$ws.Range("D14").Value = "'14.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.62%  "

$ws.Range("D15").Value = "'0.862"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.26%  "

$ws.Range("D16").Value = "2.249.51"
$ws.Range("E16").Value = "  +1.35%  "

$ws.Range("D17").Value = "42.424.15"
$ws.Range("E17").Value = "  +1.30%  "

$ws.Range("E18").Value = "  +3.92%  "

$ws.Range("D19").Value = "'6.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.53%  "

$ws.Range("D20").Value = "'72.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.95%  "

$ws.Range("D21").Value = "'11.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +54.05%  "

$ws.Range("E22").Value = "  +0.89%  "

$ws.Range("D23").Value = "'232.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.77%  "

$ws.Range("D24").Value = "'11.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.49%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("E26").Value = "  -1.68%  "

$ws.Range("D27").Value = "'2.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.29%  "

$ws.Range("E28").Value = "  +5.45%  "

$ws.Range("D29").Value = "'167.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.96%  "

$ws.Range("D30").Value = "'20.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("D31").Value = "'0.0827"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.74%  "

$ws.Range("D32").Value = "'32.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.50%  "

$ws.Range("D33").Value = "'5.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +15.62%  "

$ws.Range("E34").Value = "  +0.36%  "

$ws.Range("D35").Value = "'0.125"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.94%  "

$ws.Range("D36").Value = "'4.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.79%  "

$ws.Range("E37").Value = "  +5.83%  "

$ws.Range("D38").Value = "'14.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.59%  "

$ws.Range("E39").Value = "  +0.48%  "

$ws.Range("D40").Value = "'5.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.63%  "

$ws.Range("D41").Value = "'62.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.78%  "

$ws.Range("D42").Value = "'0.203"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").Value = "'108.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.71%  "

$ws.Range("D44").Value = "'8.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.45%  "

$ws.Range("E45").Value = "  +2.36%  "

$ws.Range("D46").Value = "'0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("D47").Value = "'2.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.49%  "

$ws.Range("E48").Value = "  +1.09%  "

$ws.Range("E49").Value = "  +2.16%  "

$ws.Range("D50").Value = "'4.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.10%  "
